$d = $word.ActiveDocument

$replacements = @(
    @("821×3=", "431×8="),
    @("714×5=", "153×3="),
    @("587×7=", "315×6="),
    @("827×9=", "716×2="),
    @("527×2=", "443×6="),
    @("928×4=", "379×7="),
    @("445×8=", "932×9="),
    @("130×4=", "982×6="),
    @("229×7=", "685×8="),
    @("807×2=", "465×7="),
    @("392×7=", "616×4="),
    @("196×6=", "754×3="),
    @("601×4=", "862×2="),
    @("131×8=", "333×8="),
    @("183×3=", "944×3="),
    @("336×2=", "222×9="),
    @("659×3=", "990×2="),
    @("548×5=", "807×4="),
    @("179×8=", "190×5="),
    @("833×9=", "757×5="),
    @("701×6=", "166×2="),
    @("648×7=", "722×6="),
    @("503×9=", "642×6="),
    @("858×9=", "932×2="),
    @("208×7=", "390×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
